$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header columns F, G, H, copying the format of the existing header cells ---
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update header row text (B1:H1) ---
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# --- Remove the now-unused 8th data row (old row 9) ---
$ws.Rows.Item(9).Delete()

# --- Row 2 (A2=0, LR) ---
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.919474060538651
$ws.Range("D2").Value = 0.01328949057939267
$ws.Range("E2").Value = 0.9075959323111146
$ws.Range("F2").Value = 0.008079828030440557
$ws.Range("G2").Value = 0.9010937210349302
$ws.Range("H2").Value = 0.01176614842018948

# --- Row 3 (A3=1, LDA) ---
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.9152424459098011
$ws.Range("D3").Value = 0.01162965171394434
$ws.Range("E3").Value = 0.8966950027806467
$ws.Range("F3").Value = 0.008827731746952868
$ws.Range("G3").Value = 0.8981608008262494
$ws.Range("H3").Value = 0.00696174371224215

# --- Row 4 (A4=2, KNN) ---
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8557011202033845
$ws.Range("D4").Value = 0.02143006205360275
$ws.Range("E4").Value = 0.8731039961865417
$ws.Range("F4").Value = 0.01157947846223474
$ws.Range("G4").Value = 0.865302825666693
$ws.Range("H4").Value = 0.008892490079090905

# --- Row 5 (A5=3, DTREE) ---
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.8914898969836604
$ws.Range("D5").Value = 0.01234489118367503
$ws.Range("E5").Value = 0.913453563200127
$ws.Range("F5").Value = 0.008843763376340269
$ws.Range("G5").Value = 0.8771761870713168
$ws.Range("H5").Value = 0.009229852154048738

# --- Row 6 (A6=4, RTREE) ---
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8719694923333599
$ws.Range("D6").Value = 0.01193368837358992
$ws.Range("E6").Value = 0.8861213950901723
$ws.Range("F6").Value = 0.008990461555954682
$ws.Range("G6").Value = 0.881408331347157
$ws.Range("H6").Value = 0.01063166626720436

# --- Row 7 (A7=5, XTREE) ---
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.9181724530600356
$ws.Range("D7").Value = 0.01407736232293906
$ws.Range("E7").Value = 0.9160541299224066
$ws.Range("F7").Value = 0.006172636980767915
$ws.Range("G7").Value = 0.9072728476470431
$ws.Range("H7").Value = 0.01118686490510426

# --- Row 8 (A8=6, SVM) ---
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.9176833240645109
$ws.Range("D8").Value = 0.01197651965327645
$ws.Range("E8").Value = 0.9253274542517411
$ws.Range("F8").Value = 0.005447372876703471
$ws.Range("G8").Value = 0.9128028918725668
$ws.Range("H8").Value = 0.01064670719590253
